$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.318.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.681.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "677.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.494"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.301.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.674.29"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.268.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.824.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").Value = "  -4.09%  "
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  -4.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.671.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "171.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.942"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000281"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.28%  "
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("E51").Value = "  -2.91%  "
